# Update gh-pages to output generated at 456a3b4
#
# The "南宁·0713国乙ONLY" event (2024-07-13) has been removed from both the
# "展览" (sheet 1) and "全部类型" (sheet 4) listings; every later row shifts
# up one position, the running index in column A is renumbered to match,
# and the refreshed "想去人数" (want-to-go count) figures from column F are
# written in for every remaining row (these are live counters that moved
# between scrapes).

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (index 1) ----------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Remove the whole row for 2024-07-13 "南宁·0713国乙ONLY"; Excel shifts
# everything below it up by one row automatically.
$ws1.Rows.Item(3).Delete()

# Renumber the running index in column A (row r now holds value r-1).
for ($r = 3; $r -le 12; $r++) {
    $ws1.Cells.Item($r, 1).Value = $r - 1
}

# Refresh the "想去人数" counts in column F.
$ws1.Cells.Item(2, 6).Value = 1595
$ws1.Cells.Item(3, 6).Value = 370
$ws1.Cells.Item(4, 6).Value = 5201
$ws1.Cells.Item(5, 6).Value = 558
$ws1.Cells.Item(6, 6).Value = 10321
$ws1.Cells.Item(7, 6).Value = 265
$ws1.Cells.Item(8, 6).Value = 562
$ws1.Cells.Item(9, 6).Value = 113
$ws1.Cells.Item(10, 6).Value = 105
$ws1.Cells.Item(11, 6).Value = 814
$ws1.Cells.Item(12, 6).Value = 82

# ---- Sheet "全部类型" (index 4) -------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# Same removal of the 2024-07-13 "南宁·0713国乙ONLY" row.
$ws4.Rows.Item(3).Delete()

# Renumber the running index in column A (row r now holds value r-1).
for ($r = 3; $r -le 18; $r++) {
    $ws4.Cells.Item($r, 1).Value = $r - 1
}

# Refresh the "想去人数" counts in column F.
$ws4.Cells.Item(2, 6).Value = 1595
$ws4.Cells.Item(3, 6).Value = 370
$ws4.Cells.Item(4, 6).Value = 14
$ws4.Cells.Item(5, 6).Value = 21
$ws4.Cells.Item(6, 6).Value = 5201
$ws4.Cells.Item(7, 6).Value = 558
$ws4.Cells.Item(8, 6).Value = 16
$ws4.Cells.Item(9, 6).Value = 10321
$ws4.Cells.Item(10, 6).Value = 265
$ws4.Cells.Item(11, 6).Value = 562
$ws4.Cells.Item(12, 6).Value = 113
$ws4.Cells.Item(13, 6).Value = 9
$ws4.Cells.Item(14, 6).Value = 3
$ws4.Cells.Item(15, 6).Value = 105
$ws4.Cells.Item(16, 6).Value = 814
$ws4.Cells.Item(17, 6).Value = 1
$ws4.Cells.Item(18, 6).Value = 82
